$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 191 (shifts existing rows 191-201 down to 192-202)
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new weekly price record
$ws.Cells.Item(191, 1).Value = 10
$ws.Cells.Item(191, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(191, 3).Value = "La Araucanía"
$ws.Cells.Item(191, 4).Value = 44585
$ws.Cells.Item(191, 5).Value = 9
$ws.Cells.Item(191, 6).Value = 100112039
$ws.Cells.Item(191, 7).Value = "Ciboulette"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 35
$ws.Cells.Item(191, 11).Value = 5000
$ws.Cells.Item(191, 12).Value = 5000
$ws.Cells.Item(191, 13).Value = 5000
$ws.Cells.Item(191, 14).Value = "`$/docena de atados"
$ws.Cells.Item(191, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(191, 16).Value = 1667
$ws.Cells.Item(191, 17).Value = 3
$ws.Cells.Item(191, 18).Value = "Hortaliza"
